# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) and "Valor Mora" (column F) values of the
# two data rows (16 and 17) are swapped:
#   Row 16: 2505 / 56940  ->  2504 / 37960
#   Row 17: 2504 / 37960  ->  2505 / 56940

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2504"
$ws.Range("F16").Value = 37960

$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 56940
